# Portugal Segunda Liga - update odds/results bases (20-06-2024)
#
# This edit corrects a data mix-up between "Academico Viseu" and "FC Porto B"
# fixtures: several whole rows had their match data (id, teams, score, odds)
# swapped with a sibling row that shares the same date, and two rows had
# their Home/Away teams swapped in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to swap for a "full row" swap: B (id) and E..AD (teams, score, odds).
# A (row index), C (Div) and D (Date) are left untouched.
$swapCols = @(2,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30)

function Swap-Rows($rowA, $rowB, $cols) {
    foreach ($c in $cols) {
        $valA = $ws.Cells.Item($rowA, $c).Value2
        $valB = $ws.Cells.Item($rowB, $c).Value2
        $ws.Cells.Item($rowA, $c).Value = $valB
        $ws.Cells.Item($rowB, $c).Value = $valA
    }
}

# Full-row swaps (same matchday pairs whose data was crossed):
Swap-Rows 5 6 $swapCols
Swap-Rows 132 133 $swapCols
Swap-Rows 140 141 $swapCols
Swap-Rows 243 244 $swapCols

# In-place Home/Away swaps (team names only needed to flip within the row):
$rows = @(99, 253)
foreach ($r in $rows) {
    $home = $ws.Cells.Item($r, 5).Value2
    $away = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 5).Value = $away
    $ws.Cells.Item($r, 6).Value = $home
}
